$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(151, "06-08-2021", 776.6900000000001),
    @(152, "09-08-2021", 784.5),
    @(153, "10-08-2021", 787.4),
    @(154, "11-08-2021", 780.08),
    @(155, "12-08-2021", 773.01),
    @(156, "13-08-2021", 772.78),
    @(157, "16-08-2021", 774.25),
    @(158, "17-08-2021", 784.1799999999999),
    @(159, "18-08-2021", 789.62),
    @(160, "19-08-2021", 787.71),
    @(161, "20-08-2021", 789.98),
    @(162, "23-08-2021", 786.4299999999999),
    @(163, "24-08-2021", 783.3),
    @(164, "25-08-2021", 782.17),
    @(165, "26-08-2021", 783.03),
    @(166, "27-08-2021", 785.0599999999999),
    @(167, "30-08-2021", 784.39),
    @(168, "31-08-2021", 779.97),
    @(169, "01-09-2021", 775.14),
    @(170, "02-09-2021", 767.1),
    @(171, "03-09-2021", 768.36),
    @(172, "06-09-2021", 766.53)
)

# Pre-format column A for the new rows as Text so the date-like strings
# ("dd-mm-yyyy") are stored as literal text, matching the source data,
# instead of being auto-converted into date serial numbers.
$dateRange = $ws.Range("A151:A172")
$dateRange.NumberFormat = "@"

foreach ($item in $data) {
    $r = $item[0]
    $dateStr = $item[1]
    $val = $item[2]

    $ws.Cells.Item($r, 1).Value = $dateStr
    $ws.Cells.Item($r, 2).Value = $val
}

# Restore the default (Normal) cell style on column A so the new cells
# do not keep a lingering explicit "Text" style attribute.
$dateRange.Style = "Normal"
